$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49..146 down to 50..147
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly data record
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D49").Value = 44973
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100112024
$ws.Range("G49").Value = "Choclo"
$ws.Range("H49").Value = "Choclero"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 2700
$ws.Range("K49").Value = 500
$ws.Range("L49").Value = 600
$ws.Range("M49").Value = 544
$ws.Range("N49").Value = "$/unidad"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 544
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
